# Fruta / hortaliza, semanal
# Update weekly price records (rows 2-7 and 10-14) on the active sheet.
# Columns: D=Fecha, M=Volumen, N=Precio minimo, O=Precio maximo,
#          P=Precio promedio ponderado, Q=Unidad de comercializacion, S=Precio $/Kg

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44210
$ws.Range("M2").Value = 70
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 11000
$ws.Range("P2").Value = 10357
$ws.Range("Q2").Value = "$/caja 14 kilos empedrada"
$ws.Range("S2").Value = 740

$ws.Range("D3").Value = 45155
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/caja 14 kilos empedrada"
$ws.Range("S3").Value = 1071

$ws.Range("D4").Value = 45142
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/caja 14 kilos empedrada"
$ws.Range("S4").Value = 1071

$ws.Range("D5").Value = 45142
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 14000
$ws.Range("Q5").Value = "$/caja 14 kilos granel"
$ws.Range("S5").Value = 1000

$ws.Range("D6").Value = 44253
$ws.Range("M6").Value = 90
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 12667
$ws.Range("Q6").Value = "$/caja 14 kilos empedrada"
$ws.Range("S6").Value = 905

$ws.Range("D7").Value = 44216
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11545
$ws.Range("Q7").Value = "$/caja 14 kilos empedrada"
$ws.Range("S7").Value = 825

$ws.Range("D10").Value = 45138
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 14000
$ws.Range("Q10").Value = "$/caja 14 kilos granel"
$ws.Range("S10").Value = 1000

$ws.Range("D11").Value = 45140
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("Q11").Value = "$/caja 14 kilos granel"
$ws.Range("S11").Value = 1071

$ws.Range("D12").Value = 44172
$ws.Range("M12").Value = 90
$ws.Range("N12").Value = 8500
$ws.Range("O12").Value = 9000
$ws.Range("P12").Value = 8806
$ws.Range("Q12").Value = "$/caja 14 kilos empedrada"
$ws.Range("S12").Value = 629

$ws.Range("D13").Value = 45152
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 16000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 16000
$ws.Range("Q13").Value = "$/caja 14 kilos empedrada"
$ws.Range("S13").Value = 1143

$ws.Range("D14").Value = 44181
$ws.Range("M14").Value = 65
$ws.Range("N14").Value = 9000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 9462
$ws.Range("Q14").Value = "$/caja 14 kilos empedrada"
$ws.Range("S14").Value = 676
